$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column widths (best effort - engine quantizes ColumnWidth input to 1/6 char
# steps, so we pick the input that lands closest to the target stored width)
# ---------------------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth  = 22         # F  -> target 22.85546875
$ws.Columns.Item(7).ColumnWidth  = 22         # G  -> target 22.85546875
$ws.Columns.Item(8).ColumnWidth  = 25.17      # H  -> target 26
$ws.Columns.Item(10).ColumnWidth = 16         # J  -> target 16.85546875 (was col F's old width)

# ---------------------------------------------------------------------------
# Row 1 - headers. E/F/G/H/J/K get rearranged: a new "send time" column is
# inserted before Destination, a new "arrival time"/"duration time" pair is
# inserted after it, and the old Events/servicePoint headers shift to J/K.
# ---------------------------------------------------------------------------
$ws.Range("E1").Value = "send time"
$ws.Range("F1").Value = "Destination"
$ws.Range("G1").Value = "arrival time"
$ws.Range("H1").Value = "duration time"
$ws.Range("J1").Value = "Events (newest)"
$ws.Range("K1").Value = "servicePoint"

# ---------------------------------------------------------------------------
# Row 2 - existing tracking number row gets the same treatment: the long
# JSON-looking event blob in F2 is replaced by send/arrival timestamps and a
# computed duration.
# ---------------------------------------------------------------------------
$ws.Range("E2").Value = "2024-02-07T17:25:00"
$ws.Range("F2").Value = "MUNICH - GERMANY"
$ws.Range("G2").Value = "2024-02-08T08:28:00"
$ws.Range("H2").Value = 0.6270833333333333
$ws.Range("H2").NumberFormat = "[hh]:mm:ss"

# ---------------------------------------------------------------------------
# Rows 3-6 - new tracking rows with send/arrival/duration data.
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = 8585095551
$ws.Range("C3").Value = "delivered"
$ws.Range("D3").Value = "MUNICH - GERMANY"
$ws.Range("E3").Value = "2024-02-26T19:09:00"
$ws.Range("F3").Value = "CAMBRIDGE - UK"
$ws.Range("G3").Value = "2024-02-27T08:58:00"
$ws.Range("H3").Value = 0.5756944444444444
$ws.Range("H3").NumberFormat = "[hh]:mm:ss"

$ws.Range("A4").Value = 4874299916
$ws.Range("C4").Value = "delivered"
$ws.Range("D4").Value = "MUNICH - GERMANY"
$ws.Range("E4").Value = "2024-02-26T16:20:00"
$ws.Range("F4").Value = "VERONA - ITALY"
$ws.Range("G4").Value = "2024-02-27T12:14:00"
$ws.Range("H4").Value = 0.8291666666666667
$ws.Range("H4").NumberFormat = "[hh]:mm:ss"

$ws.Range("A5").Value = 3527880731
$ws.Range("C5").Value = "delivered"
$ws.Range("D5").Value = "MUNICH - GERMANY"
$ws.Range("E5").Value = "2024-02-28T17:56:00"
$ws.Range("F5").Value = "CAMBRIDGE - UK"
$ws.Range("G5").Value = "2024-03-01T08:47:00"
$ws.Range("H5").Value = 1.61875
$ws.Range("H5").NumberFormat = "[hh]:mm:ss"

$ws.Range("A6").Value = 7177497666
$ws.Range("C6").Value = "delivered"
$ws.Range("D6").Value = "MUNICH - GERMANY"
$ws.Range("E6").Value = 45350.68472222222
$ws.Range("E6").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("F6").Value = "VERONA - ITALY"
$ws.Range("G6").Value = 45351.43194444444
$ws.Range("G6").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("H6").Value = 0.7472222222222222
$ws.Range("H6").NumberFormat = "[hh]:mm:ss"

# ---------------------------------------------------------------------------
# Rows 7-11 - more tracking numbers queued up (no further details yet).
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = 89352111143
$ws.Range("A8").Value = 6113737884
$ws.Range("A9").Value = 6691061322
$ws.Range("A10").Value = 2685659756
$ws.Range("A11").Value = 6691061322

# ---------------------------------------------------------------------------
# Old example block (previously A15:G16) moves out to N14:T15, well clear of
# the real data table, and the old "A16 / id" example cell is replaced by an
# actual tracking number.
# ---------------------------------------------------------------------------
$ws.Range("A15").ClearContents()
$ws.Range("B16").ClearContents()
$ws.Range("C16").ClearContents()
$ws.Range("D16").ClearContents()
$ws.Range("E16").ClearContents()
$ws.Range("F16").ClearContents()
$ws.Range("G16").ClearContents()

$ws.Range("N14").Value = "E.g"

$ws.Range("N15").Value = '"id":"6921406671"'
$ws.Range("O15").Value = '"service":"express"'
$ws.Range("P15").Value = '"status"{status: "delivered"} '
$ws.Range("Q15").Value = '"origin": {address:{"addresslocality": XXXX - SPAIN'
$ws.Range("R15").Value = '"Destination": {address:{"addresslocality": XXXX - SPAIN'
$ws.Range("S15").Value = '(the newest updata)  events[{0}]'
$ws.Range("T15").Value = '"servicePoint":{"url":"http://www.dhl.es/en/country_profile.html"}'

$ws.Range("A16").Value = 11111111111

# ---------------------------------------------------------------------------
# Final selection, to match the saved view state.
# ---------------------------------------------------------------------------
$ws.Range("A12").Select() | Out-Null
